$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values look like plain numbers need to be forced to text
# so Excel does not silently convert them to numeric cells (losing formatting
# such as trailing zeros), matching the inline-string cells used in the source file.
$textCells = @("D5","D6","D10","D13","D14","D17","D19","D20","D22","D23","D27","D28","D29","D31","D33","D34","D37","D40","D44","D45","D46","D47","D48","D49")
foreach ($ref in $textCells) {
  $ws.Range($ref).NumberFormat = "@"
}

# Apply the updated crypto market data
$ws.Range("D2").Value = "43.193.83"
$ws.Range("E2").Value = "  +0.29%  "

$ws.Range("D3").Value = "2.322.40"
$ws.Range("E3").Value = "  +0.79%  "

$ws.Range("E4").Value = "  +0.00%  "

$ws.Range("D5").Value = "302.55"
$ws.Range("E5").Value = "  +0.06%  "

$ws.Range("D6").Value = "99.44"
$ws.Range("E6").Value = "  +0.02%  "

$ws.Range("E7").Value = "  +0.12%  "

$ws.Range("E9").Value = "  +1.75%  "

$ws.Range("D10").Value = "36.27"
$ws.Range("E10").Value = "  +5.38%  "

$ws.Range("E11").Value = "  -0.72%  "

$ws.Range("E12").Value = "  -0.91%  "

$ws.Range("D13").Value = "17.64"
$ws.Range("E13").Value = "  -0.80%  "

$ws.Range("D14").Value = "6.93"
$ws.Range("E14").Value = "  +1.86%  "

$ws.Range("D15").Value = "2.683.33"
$ws.Range("E15").Value = "  +0.77%  "

$ws.Range("D16").Value = "2.273.38"
$ws.Range("E16").Value = "  -0.85%  "

$ws.Range("D17").Value = "0.798"
$ws.Range("E17").Value = "  -1.26%  "

$ws.Range("D18").Value = "43.108.29"
$ws.Range("E18").Value = "  +0.37%  "

$ws.Range("D19").Value = "12.96"
$ws.Range("E19").Value = "  +4.85%  "

$ws.Range("D20").Value = "6.23"
$ws.Range("E20").Value = "  +2.00%  "

$ws.Range("E21").Value = "  +0.40%  "

$ws.Range("D22").Value = "68.17"
$ws.Range("E22").Value = "  +0.52%  "

$ws.Range("D23").Value = "240.64"
$ws.Range("E23").Value = "  +1.53%  "

$ws.Range("E24").Value = "  -0.77%  "

$ws.Range("E25").Value = "  -0.89%  "

$ws.Range("E26").Value = "  -0.08%  "

$ws.Range("D27").Value = "25.55"
$ws.Range("E27").Value = "  +4.06%  "

$ws.Range("D28").Value = "168.99"
$ws.Range("E28").Value = "  +0.40%  "

$ws.Range("D29").Value = "34.19"
$ws.Range("E29").Value = "  +0.96%  "

$ws.Range("E30").Value = "  +0.26%  "

$ws.Range("D31").Value = "2.04"
$ws.Range("E31").Value = "  -2.51%  "

$ws.Range("E32").Value = "  +3.42%  "

$ws.Range("D33").Value = "0.999"
$ws.Range("E33").Value = "  -0.09%  "

$ws.Range("D34").Value = "4.73"
$ws.Range("E34").Value = "  +3.69%  "

$ws.Range("E35").Value = "  +4.25%  "

$ws.Range("E36").Value = "  -0.92%  "

$ws.Range("D37").Value = "0.0697"
$ws.Range("E37").Value = "  -0.45%  "

$ws.Range("E38").Value = "  +0.32%  "

$ws.Range("E39").Value = "  +0.54%  "

$ws.Range("D40").Value = "2.77"
$ws.Range("E40").Value = "  -1.74%  "

$ws.Range("E41").Value = "  +0.33%  "

$ws.Range("D42").Value = "2.002.16"
$ws.Range("E42").Value = "  -0.02%  "

$ws.Range("E43").Value = "  +1.50%  "

$ws.Range("D44").Value = "2.24"
$ws.Range("E44").Value = "  -4.48%  "

$ws.Range("D45").Value = "10.10"

$ws.Range("D46").Value = "17.62"
$ws.Range("E46").Value = "  -0.34%  "

$ws.Range("D47").Value = "2.87"
$ws.Range("E47").Value = "  +0.29%  "

$ws.Range("B48").Value = "BitcoinSV"
$ws.Range("C48").Value = "https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv"
$ws.Range("D48").Value = "76.66"
$ws.Range("E48").Value = "  +9.15%  "

$ws.Range("B49").Value = "MultiversX"
$ws.Range("C49").Value = "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
$ws.Range("D49").Value = "54.96"
$ws.Range("E49").Value = "  -0.98%  "

$ws.Range("D50").Value = "2.548.27"
$ws.Range("E50").Value = "  +0.83%  "

$ws.Range("E51").Value = "  +1.87%  "

# Restore default (General) style bookkeeping so no stray number format remains
foreach ($ref in $textCells) {
  $ws.Range($ref).Style = "Normal"
}
